$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6658955216407776
$ws.Range("B1").Value = 3.069838285446167
$ws.Range("C1").Value = 4.542080879211426
$ws.Range("D1").Value = 2.374843835830688
$ws.Range("E1").Value = 1.063945770263672
